$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(8, 9).Value = 'sd'
$ws.Cells.Item(8, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(14, 9).Value = 'sd'
$ws.Cells.Item(14, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(22, 9).Value = 'sd'
$ws.Cells.Item(22, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(23, 9).Value = '%'
$ws.Cells.Item(23, 10).Value = 'Uninterpretable'
$ws.Cells.Item(24, 9).Value = 'b'
$ws.Cells.Item(24, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(25, 9).Value = 'sd'
$ws.Cells.Item(25, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(47, 9).Value = 'sd'
$ws.Cells.Item(47, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(48, 9).Value = 'b'
$ws.Cells.Item(48, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(80, 9).Value = 'sd'
$ws.Cells.Item(80, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(83, 9).Value = 'sd'
$ws.Cells.Item(83, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(85, 9).Value = 'aa'
$ws.Cells.Item(85, 10).Value = 'Agree/Accept'
$ws.Cells.Item(99, 9).Value = 'sd'
$ws.Cells.Item(99, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(112, 9).Value = 'sd'
$ws.Cells.Item(112, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(126, 9).Value = 'aa'
$ws.Cells.Item(126, 10).Value = 'Agree/Accept'
$ws.Cells.Item(127, 9).Value = 'sd'
$ws.Cells.Item(127, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(130, 9).Value = '%'
$ws.Cells.Item(130, 10).Value = 'Uninterpretable'
$ws.Cells.Item(134, 9).Value = 'ba'
$ws.Cells.Item(134, 10).Value = 'Appreciation'
$ws.Cells.Item(135, 9).Value = 'b'
$ws.Cells.Item(135, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(138, 9).Value = 'aa'
$ws.Cells.Item(138, 10).Value = 'Agree/Accept'
$ws.Cells.Item(139, 9).Value = 'sd'
$ws.Cells.Item(139, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(141, 9).Value = 'aa'
$ws.Cells.Item(141, 10).Value = 'Agree/Accept'
$ws.Cells.Item(147, 9).Value = '%'
$ws.Cells.Item(147, 10).Value = 'Uninterpretable'
$ws.Cells.Item(148, 9).Value = 'sd'
$ws.Cells.Item(148, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(153, 9).Value = 'b'
$ws.Cells.Item(153, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(174, 9).Value = 'ba'
$ws.Cells.Item(174, 10).Value = 'Appreciation'
$ws.Cells.Item(188, 9).Value = 'sd'
$ws.Cells.Item(188, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(190, 9).Value = 'sv'
$ws.Cells.Item(190, 10).Value = 'Statement-opinion'
$ws.Cells.Item(204, 9).Value = 'sv'
$ws.Cells.Item(204, 10).Value = 'Statement-opinion'
$ws.Cells.Item(205, 9).Value = 'b'
$ws.Cells.Item(205, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(210, 9).Value = 'b'
$ws.Cells.Item(210, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(216, 9).Value = 'ba'
$ws.Cells.Item(216, 10).Value = 'Appreciation'
$ws.Cells.Item(217, 9).Value = 'sv'
$ws.Cells.Item(217, 10).Value = 'Statement-opinion'
$ws.Cells.Item(227, 9).Value = 'aa'
$ws.Cells.Item(227, 10).Value = 'Agree/Accept'
$ws.Cells.Item(230, 9).Value = 'ba'
$ws.Cells.Item(230, 10).Value = 'Appreciation'
$ws.Cells.Item(231, 9).Value = '%'
$ws.Cells.Item(231, 10).Value = 'Uninterpretable'
$ws.Cells.Item(235, 9).Value = 'sd'
$ws.Cells.Item(235, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(241, 9).Value = 'sv'
$ws.Cells.Item(241, 10).Value = 'Statement-opinion'
$ws.Cells.Item(245, 9).Value = 'sd'
$ws.Cells.Item(245, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(254, 9).Value = 'sv'
$ws.Cells.Item(254, 10).Value = 'Statement-opinion'
$ws.Cells.Item(259, 9).Value = 'ba'
$ws.Cells.Item(259, 10).Value = 'Appreciation'
$ws.Cells.Item(284, 9).Value = 'sv'
$ws.Cells.Item(284, 10).Value = 'Statement-opinion'
$ws.Cells.Item(287, 9).Value = 'b'
$ws.Cells.Item(287, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(296, 9).Value = 'ba'
$ws.Cells.Item(296, 10).Value = 'Appreciation'
$ws.Cells.Item(299, 9).Value = 'ba'
$ws.Cells.Item(299, 10).Value = 'Appreciation'
$ws.Cells.Item(308, 9).Value = 'ba'
$ws.Cells.Item(308, 10).Value = 'Appreciation'
$ws.Cells.Item(311, 9).Value = 'sd'
$ws.Cells.Item(311, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(317, 9).Value = 'sd'
$ws.Cells.Item(317, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(321, 9).Value = 'sd'
$ws.Cells.Item(321, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(322, 9).Value = 'ba'
$ws.Cells.Item(322, 10).Value = 'Appreciation'
$ws.Cells.Item(328, 9).Value = 'sd'
$ws.Cells.Item(328, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(347, 9).Value = 'sd'
$ws.Cells.Item(347, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(349, 9).Value = 'b'
$ws.Cells.Item(349, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(352, 9).Value = 'sd'
$ws.Cells.Item(352, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(379, 9).Value = 'sd'
$ws.Cells.Item(379, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(384, 9).Value = 'ba'
$ws.Cells.Item(384, 10).Value = 'Appreciation'
$ws.Cells.Item(404, 9).Value = 'sv'
$ws.Cells.Item(404, 10).Value = 'Statement-opinion'
$ws.Cells.Item(406, 9).Value = 'ba'
$ws.Cells.Item(406, 10).Value = 'Appreciation'
$ws.Cells.Item(407, 9).Value = 'b'
$ws.Cells.Item(407, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(419, 9).Value = 'b'
$ws.Cells.Item(419, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(420, 9).Value = 'sv'
$ws.Cells.Item(420, 10).Value = 'Statement-opinion'
$ws.Cells.Item(423, 9).Value = 'sd'
$ws.Cells.Item(423, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(426, 9).Value = 'aa'
$ws.Cells.Item(426, 10).Value = 'Agree/Accept'
$ws.Cells.Item(429, 9).Value = 'sd'
$ws.Cells.Item(429, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(440, 9).Value = 'sv'
$ws.Cells.Item(440, 10).Value = 'Statement-opinion'
$ws.Cells.Item(449, 9).Value = 'b'
$ws.Cells.Item(449, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(452, 9).Value = 'sd'
$ws.Cells.Item(452, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(456, 9).Value = 'sd'
$ws.Cells.Item(456, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(466, 9).Value = 'sd'
$ws.Cells.Item(466, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(468, 9).Value = 'sd'
$ws.Cells.Item(468, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(470, 9).Value = 'sd'
$ws.Cells.Item(470, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(476, 9).Value = 'sd'
$ws.Cells.Item(476, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(482, 9).Value = 'sv'
$ws.Cells.Item(482, 10).Value = 'Statement-opinion'
$ws.Cells.Item(487, 9).Value = 'sd'
$ws.Cells.Item(487, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(492, 9).Value = 'sv'
$ws.Cells.Item(492, 10).Value = 'Statement-opinion'
$ws.Cells.Item(493, 9).Value = 'aa'
$ws.Cells.Item(493, 10).Value = 'Agree/Accept'
$ws.Cells.Item(499, 9).Value = 'sv'
$ws.Cells.Item(499, 10).Value = 'Statement-opinion'
$ws.Cells.Item(507, 9).Value = 'sd'
$ws.Cells.Item(507, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(522, 9).Value = 'sv'
$ws.Cells.Item(522, 10).Value = 'Statement-opinion'
$ws.Cells.Item(530, 9).Value = 'sv'
$ws.Cells.Item(530, 10).Value = 'Statement-opinion'
$ws.Cells.Item(531, 9).Value = 'sd'
$ws.Cells.Item(531, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(533, 9).Value = 'aa'
$ws.Cells.Item(533, 10).Value = 'Agree/Accept'
$ws.Cells.Item(548, 9).Value = '%'
$ws.Cells.Item(548, 10).Value = 'Uninterpretable'
$ws.Cells.Item(552, 9).Value = 'sv'
$ws.Cells.Item(552, 10).Value = 'Statement-opinion'
$ws.Cells.Item(555, 9).Value = 'ba'
$ws.Cells.Item(555, 10).Value = 'Appreciation'
